$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are written as literal text, matching the
# source data (which stores prices like "64.121.45" / "525.85" as strings,
# not numbers) -- without this Excel auto-converts numeric-looking text.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.121.45"
$ws.Range("E2").Value = "  -1.31%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.300.31"
$ws.Range("E3").Value = "  -0.58%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.85"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.13"
$ws.Range("E6").Value = "  -5.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.584"
$ws.Range("E7").Value = "  -3.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.294.93"
$ws.Range("E8").Value = "  -0.62%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.601"
$ws.Range("E10").Value = "  -3.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.69"
$ws.Range("E11").Value = "  -12.05%  "

$ws.Range("E12").Value = "  -1.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("E13").Value = "  -2.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.89"
$ws.Range("E14").Value = "  -2.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.836.76"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.309.40"
$ws.Range("E16").Value = "  -0.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.117"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.079.67"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.31"
$ws.Range("E19").Value = "  -2.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.09"
$ws.Range("E20").Value = "  -1.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.949"
$ws.Range("E21").Value = "  -1.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.79"
$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("E23").Value = "  +6.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.02"
$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.69"
$ws.Range("E25").Value = "  -3.93%  "

$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.03"
$ws.Range("E26").Value = "  -1.91%  "

$ws.Range("E27").Value = "  +1.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.69"
$ws.Range("E28").Value = "  -0.52%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.17"
$ws.Range("E29").Value = "  -4.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.09"
$ws.Range("E30").Value = "  -4.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.64"
$ws.Range("E31").Value = "  -1.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "626.42"
$ws.Range("E32").Value = "  -4.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.52"
$ws.Range("E33").Value = "  -4.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.10"
$ws.Range("E34").Value = "  -2.43%  "

$ws.Range("E35").Value = "  -1.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.01"
$ws.Range("E36").Value = "  -4.22%  "

$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.81"
$ws.Range("E38").Value = "  -2.56%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.376"
$ws.Range("E39").Value = "  -4.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0746"
$ws.Range("E40").Value = "  +4.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.19"
$ws.Range("E42").Value = "  +9.61%  "

$ws.Range("E43").Value = "  -2.41%  "

$ws.Range("E44").Value = "  +3.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.888.20"
$ws.Range("E45").Value = "  -0.40%  "

$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.13"
$ws.Range("E47").Value = "  +4.83%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0395"
$ws.Range("E48").Value = "  -1.71%  "

$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.61"
$ws.Range("E49").Value = "  -4.27%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "137.84"
$ws.Range("E50").Value = "  +1.65%  "

$ws.Range("E51").Value = "  -2.31%  "
